$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.979.16'
$ws.Range('E2').Value = '  -2.01%  '
$ws.Range('D3').Value = '1.650.03'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.69'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3895'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3818'
$ws.Range('E8').Value = '  -2.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.33'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('E10').Value = '  -4.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.001'
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08447'
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.78'
$ws.Range('E13').Value = '  -3.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.074'
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.980'
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001310'
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('D17').Value = '1.656.86'
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.57'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06984'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.64'
$ws.Range('E20').Value = '  -4.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.974'
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.78'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('D24').Value = '23.979.37'
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.444'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.974'
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('E27').Value = '  -2.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.19'
$ws.Range('E28').Value = '  -3.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.434'
$ws.Range('E29').Value = '  +1.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '138.32'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.890'
$ws.Range('E31').Value = '  -3.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.516'
$ws.Range('E32').Value = '  -1.51%  '
$ws.Range('D33').Value = '1.830.14'
$ws.Range('E33').Value = '  -0.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.016'
$ws.Range('E34').Value = '  -4.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08079'
$ws.Range('E35').Value = '  -2.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.730'
$ws.Range('E36').Value = '  -1.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02920'
$ws.Range('E37').Value = '  -1.92%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2673'
$ws.Range('E38').Value = '  -2.88%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '10.69'
$ws.Range('E39').Value = '  -4.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09112'
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7580'
$ws.Range('E41').Value = '  -2.29%  '
$ws.Range('E42').Value = '  -3.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.420'
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.38'
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('E45').Value = '  -2.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.462'
$ws.Range('E46').Value = '  -2.81%  '
$ws.Range('E47').Value = '  -1.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08328'
$ws.Range('E49').Value = '  -1.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '134.79'
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.218'
$ws.Range('E51').Value = '  -3.71%  '
